$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strip the stray "_GoBack" bookmark out of the empty second paragraph.
#    (It gets re-added further below, inside the new "O portal..." text.)
#    Deleting the Bookmark object (rather than rewriting the paragraph's
#    XML) leaves the paragraph mark's own identity attributes untouched.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2) Append the two new paragraphs at the very end of the document body
#    (before the final sectPr), describing comments and the investment
#    simulator.
# ---------------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr>'
$pPr = '<w:pPr><w:spacing w:before="240" w:after="240" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/>' + $rPr + '</w:pPr>'

function New-Run([string]$text, [bool]$preserve) {
    if ($preserve) {
        return '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r>'
    }
}

$bookmark = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$p1 = '<w:p>' + $pPr `
    + (New-Run 'Toda notícia ' $true) `
    + (New-Run 'contém comentários' $false) `
    + (New-Run ', ' $true) `
    + (New-Run 'curtidas' $false) `
    + (New-Run ',' $false) `
    + (New-Run ' link que redireciona para ' $true) `
    + (New-Run 'ela e um ID.' $false) `
    + (New-Run ' ' $true) `
    + (New-Run 'Todo comentário tem uma data e hora, curtidas, quem ' $true) `
    + (New-Run 'e o que comentou, além das respostas que também são outros comentários.' $false) `
    + '</w:p>'

$p2 = '<w:p>' + $pPr `
    + (New-Run 'O porta' $false) `
    + (New-Run 'l' $false) `
    + $bookmark `
    + (New-Run ' também possibilita que o usuário faça simulações d' $true) `
    + (New-Run 'e investimentos, nele deve ser armazenado o valor inicial' $false) `
    + (New-Run ',' $false) `
    + (New-Run ' atual' $true) `
    + (New-Run ' (sendo este o resultado das variações durante um determinado período)' $true) `
    + (New-Run ' e data da simulação.' $true) `
    + '</w:p>'

$endRange = $d.Content
$endRange.Collapse(0)
$xmlNewParas = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $p1 + $p2 + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($xmlNewParas)
